$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "type"
$ws.Range("B6").Value = "non-ingredient"
$ws.Range("A6").Value = "paper towels"
$ws.Range("A7").Value = "toilet paper"
$ws.Range("B7").Value = "non-ingredient"

$ws.Columns.Item(2).ColumnWidth = 18.6

$ws.Range("B8").Select() | Out-Null
